$d = $word.ActiveDocument

# 1. Title: A-E-C-P -> A-E-K-P (first heading)
$d.Content.Find.Execute("A-E-C-P: Aseptá, Eksplorá, Konektá i Praktiká", $true, $false, $false, $false, $false,
                         $true, 1, $false, "A-E-K-P: Aseptá, Eksplorá, Konektá i Praktiká", 2)

# 2. Remove comma before "i permití" in the method description paragraph
$d.Content.Find.Execute("enfatisá prinsipionan sentral, i permití oportunidatnan", $true, $false, $false, $false, $false,
                         $true, 1, $false, "enfatisá prinsipionan sentral i permití oportunidatnan", 2)

# 3. "A-E-C-P ta para pa e 4 pasonan klave" -> "A-E-K-P ta para pa e 4 pasonan klave"
$d.Content.Find.Execute("A-E-C-P ta para pa e 4 pasonan klave", $true, $false, $false, $false, $false,
                         $true, 1, $false, "A-E-K-P ta para pa e 4 pasonan klave", 2)

# 4a. "A-C-E-P na kualke momento" -> "A-K-E-P na kualke momentu"
$d.Content.Find.Execute("A-C-E-P na kualke momento", $true, $false, $false, $false, $false,
                         $true, 1, $false, "A-K-E-P na kualke momentu", 2)

# 4b. "duna ehèmpel di e A-E-C-P pa asina" -> "duna ehèmpel di e A-E-K-P pa asina"
$d.Content.Find.Execute("duna ehèmpel di e A-E-C-P pa asina", $true, $false, $false, $false, $false,
                         $true, 1, $false, "duna ehèmpel di e A-E-K-P pa asina", 2)

# 5. "Thank you for sharing." -> "Danki pa kompartí."
$d.Content.Find.Execute("sigui su sugerensianan. Thank you for sharing. ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "sigui su sugerensianan. Danki pa kompartí. ", 2)

# 6. "Kon siguimentu di e ehèmpel di bo yu a laga bo sinti?" -> "Kon siguimentu di bo yu su guia a laga bo sinti?"
$d.Content.Find.Execute("Kon siguimentu di e ehèmpel di bo yu a laga bo sinti?", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Kon siguimentu di bo yu su guia a laga bo sinti?", 2)

# 7. "Enkurashá mayornan pa traha nan mes konekshonnan" -> "Enkurashá mayornan pa krea nan mes konekshonnan"
$d.Content.Find.Execute("Enkurashá mayornan pa traha nan mes konekshonnan", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Enkurashá mayornan pa krea nan mes konekshonnan", 2)

# 8. "nos atenshon kompleto." -> "nos atenshon kompletu."
$d.Content.Find.Execute("nos atenshon kompleto. Esaki ta bai bèk", $true, $false, $false, $false, $false,
                         $true, 1, $false, "nos atenshon kompletu. Esaki ta bai bèk", 2)
